# Colombia Primera A - base update (11-06-2024 21:19)
#
# The scraped source data got a handful of fixture rows reshuffled
# (their id/teams/odds moved to a different row within the same block of
# matches) while the running index in column A (and the constant Div /
# Date columns inside each block) stayed put. This script re-reads each
# affected row's B:AD payload into memory first, then writes the rows
# back out in the permuted order so every block ends up matching the
# target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B" + $row + ":AD" + $row).Value()
}

function Set-RowData($row, $data) {
    $ws.Range("B" + $row + ":AD" + $row).Value = $data
}

function Apply-Cycle($rows) {
    # $rows is an ordered list of row numbers; the content that was in
    # $rows[1] before the edit ends up in $rows[0], $rows[2]'s content
    # ends up in $rows[1], ... and $rows[0]'s original content wraps
    # around into the last row.
    $n = $rows.Length
    $snapshot = @()
    for ($i = 0; $i -lt $n; $i++) {
        $snapshot += ,(Get-RowData $rows[$i])
    }
    for ($i = 0; $i -lt $n; $i++) {
        $srcIndex = ($i + 1) % $n
        Set-RowData $rows[$i] $snapshot[$srcIndex]
    }
}

# Rows 14 & 15 swap their data (simple 2-cycle).
Apply-Cycle @(14, 15)

# Rows 208, 209, 210, 211, 212 rotate: 208<-211<-209<-210<-212<-208.
Apply-Cycle @(208, 211, 209, 210, 212)

# Rows 214 & 215 swap.
Apply-Cycle @(214, 215)

# Rows 240 & 241 swap.
Apply-Cycle @(240, 241)

# Rows 427, 428, 429, 430 rotate: 427<-428<-429<-430<-427.
Apply-Cycle @(427, 428, 429, 430)
